{"js": "// The document contains a single table of two-digit \u00f7 one-digit division\n// problems, laid out 5-per-row with blank spacer rows in between. The\n// content-bearing table rows are at indices 0, 4, 8, 12, 16 (5 cells each).\n// We update each cell's expression text in place, by (row, col) position,\n// so that duplicate expressions (e.g. \"63\u00f78=\" appears twice in the source)\n// are each mapped to their own, independent replacement.\nconst edits = [\n  { row: 0, col: 0, oldText: \"79\u00f73=\", newText: \"53\u00f75=\" },\n  { row: 0, col: 1, oldText: \"25\u00f79=\", newText: \"16\u00f73=\" },\n  { row: 0, col: 2, oldText: \"11\u00f77=\", newText: \"52\u00f73=\" },\n  { row: 0, col: 3, oldText: \"95\u00f79=\", newText: \"10\u00f75=\" },\n  { row: 0, col: 4, oldText: \"28\u00f79=\", newText: \"77\u00f76=\" },\n  { row: 4, col: 0, oldText: \"63\u00f72=\", newText: \"79\u00f76=\" },\n  { row: 4, col: 1, oldText: \"33\u00f78=\", newText: \"99\u00f77=\" },\n  { row: 4, col: 2, oldText: \"30\u00f76=\", newText: \"83\u00f76=\" },\n  { row: 4, col: 3, oldText: \"60\u00f77=\", newText: \"30\u00f79=\" },\n  { row: 4, col: 4, oldText: \"63\u00f78=\", newText: \"42\u00f78=\" },\n  { row: 8, col: 0, oldText: \"16\u00f76=\", newText: \"62\u00f79=\" },\n  { row: 8, col: 1, oldText: \"21\u00f74=\", newText: \"27\u00f76=\" },\n  { row: 8, col: 2, oldText: \"62\u00f77=\", newText: \"17\u00f73=\" },\n  { row: 8, col: 3, oldText: \"24\u00f78=\", newText: \"56\u00f75=\" },\n  { row: 8, col: 4, oldText: \"33\u00f77=\", newText: \"39\u00f75=\" },\n  { row: 12, col: 0, oldText: \"50\u00f76=\", newText: \"61\u00f76=\" },\n  { row: 12, col: 1, oldText: \"46\u00f72=\", newText: \"99\u00f74=\" },\n  { row: 12, col: 2, oldText: \"63\u00f78=\", newText: \"86\u00f76=\" },\n  { row: 12, col: 3, oldText: \"18\u00f78=\", newText: \"26\u00f75=\" },\n  { row: 12, col: 4, oldText: \"76\u00f76=\", newText: \"62\u00f79=\" },\n  { row: 16, col: 0, oldText: \"77\u00f76=\", newText: \"47\u00f79=\" },\n  { row: 16, col: 1, oldText: \"17\u00f79=\", newText: \"68\u00f77=\" },\n  { row: 16, col: 2, oldText: \"24\u00f77=\", newText: \"54\u00f73=\" },\n  { row: 16, col: 3, oldText: \"59\u00f73=\", newText: \"77\u00f79=\" },\n  { row: 16, col: 4, oldText: \"70\u00f75=\", newText: \"96\u00f74=\" },\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Sanity-check each cell still holds the expected original text before\n// overwriting it, then write the new expression straight into the cell's\n// text (this preserves the existing run/paragraph formatting, since\n// TableCell.value maps to a Range.Text-style write on the existing run).\nfor (const edit of edits) {\n  const cell = table.getCell(edit.row, edit.col);\n  cell.load(\"value\");\n  await context.sync();\n\n  if (cell.value !== edit.oldText) {\n    throw new Error(\n      `Unexpected text at row ${edit.row}, col ${edit.col}: ` +\n      `expected \"${edit.oldText}\" but found \"${cell.value}\"`\n    );\n  }\n\n  cell.value = edit.newText;\n}\n\nawait context.sync();\n", "ps1": "# The document contains a single table of two-digit \u00f7 one-digit division\n# problems, laid out 5-per-row with blank spacer rows in between. The\n# content-bearing table rows are Word rows 1, 5, 9, 13, 17 (5 cells each,\n# 1-based COM indexing). We update each cell's expression text in place by\n# (row, col) position, so that duplicate expressions (e.g. \"63\u00f78=\" appears\n# twice in the source) are each mapped to their own, independent\n# replacement.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$edits = @(\n    @{ Row = 1; Col = 1; OldText = \"79\u00f73=\"; NewText = \"53\u00f75=\" },\n    @{ Row = 1; Col = 2; OldText = \"25\u00f79=\"; NewText = \"16\u00f73=\" },\n    @{ Row = 1; Col = 3; OldText = \"11\u00f77=\"; NewText = \"52\u00f73=\" },\n    @{ Row = 1; Col = 4; OldText = \"95\u00f79=\"; NewText = \"10\u00f75=\" },\n    @{ Row = 1; Col = 5; OldText = \"28\u00f79=\"; NewText = \"77\u00f76=\" },\n    @{ Row = 5; Col = 1; OldText = \"63\u00f72=\"; NewText = \"79\u00f76=\" },\n    @{ Row = 5; Col = 2; OldText = \"33\u00f78=\"; NewText = \"99\u00f77=\" },\n    @{ Row = 5; Col = 3; OldText = \"30\u00f76=\"; NewText = \"83\u00f76=\" },\n    @{ Row = 5; Col = 4; OldText = \"60\u00f77=\"; NewText = \"30\u00f79=\" },\n    @{ Row = 5; Col = 5; OldText = \"63\u00f78=\"; NewText = \"42\u00f78=\" },\n    @{ Row = 9; Col = 1; OldText = \"16\u00f76=\"; NewText = \"62\u00f79=\" },\n    @{ Row = 9; Col = 2; OldText = \"21\u00f74=\"; NewText = \"27\u00f76=\" },\n    @{ Row = 9; Col = 3; OldText = \"62\u00f77=\"; NewText = \"17\u00f73=\" },\n    @{ Row = 9; Col = 4; OldText = \"24\u00f78=\"; NewText = \"56\u00f75=\" },\n    @{ Row = 9; Col = 5; OldText = \"33\u00f77=\"; NewText = \"39\u00f75=\" },\n    @{ Row = 13; Col = 1; OldText = \"50\u00f76=\"; NewText = \"61\u00f76=\" },\n    @{ Row = 13; Col = 2; OldText = \"46\u00f72=\"; NewText = \"99\u00f74=\" },\n    @{ Row = 13; Col = 3; OldText = \"63\u00f78=\"; NewText = \"86\u00f76=\" },\n    @{ Row = 13; Col = 4; OldText = \"18\u00f78=\"; NewText = \"26\u00f75=\" },\n    @{ Row = 13; Col = 5; OldText = \"76\u00f76=\"; NewText = \"62\u00f79=\" },\n    @{ Row = 17; Col = 1; OldText = \"77\u00f76=\"; NewText = \"47\u00f79=\" },\n    @{ Row = 17; Col = 2; OldText = \"17\u00f79=\"; NewText = \"68\u00f77=\" },\n    @{ Row = 17; Col = 3; OldText = \"24\u00f77=\"; NewText = \"54\u00f73=\" },\n    @{ Row = 17; Col = 4; OldText = \"59\u00f73=\"; NewText = \"77\u00f79=\" },\n    @{ Row = 17; Col = 5; OldText = \"70\u00f75=\"; NewText = \"96\u00f74=\" }\n)\n\nforeach ($edit in $edits) {\n    $cell = $t.Cell($edit.Row, $edit.Col)\n    # Cell.Range.Text includes the trailing cell-end mark (CR + BEL); strip\n    # it before comparing against the plain expression text.\n    $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n    if ($current -ne $edit.OldText) {\n        throw \"Unexpected text at row $($edit.Row), col $($edit.Col): expected '$($edit.OldText)' but found '$current'\"\n    }\n    # Assigning Range.Text rewrites the run's text in place, preserving the\n    # existing run/paragraph formatting (font, size, alignment, etc.).\n    $cell.Range.Text = $edit.NewText\n}\n"}
